# WS_holdings.xlsx update
# - bump the "as of" date in the confidential disclosure note (A16) from
#   2021-03-22 to 2021-03-23
# - refresh the Weight (D) / Percent Change (E) figures for rows 2-13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no-op cell editing guard); lift it so the
# cells below can be written, then restore protection afterwards.
$ws.Unprotect()

$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.03001547701628224
$ws.Range("E2").Value = -0.002214839424141579

$ws.Range("D3").Value = 0.02342230791792799
$ws.Range("E3").Value = 0.01274044466650004

$ws.Range("D4").Value = 0.05170757601938909
$ws.Range("E4").Value = 0.002455795677799744

$ws.Range("D5").Value = 0.1387049153469954
$ws.Range("E5").Value = -0.01475826972010197

$ws.Range("D6").Value = 0.03061354318600456
$ws.Range("E6").Value = -0.01901140684410652

$ws.Range("D7").Value = 0.1206511352052821
$ws.Range("E7").Value = -0.01362862010221488

$ws.Range("D8").Value = 0.100415110706978
$ws.Range("E8").Value = -0.01998041136141038

$ws.Range("D9").Value = 0.02765060087971998
$ws.Range("E9").Value = -0.02502870264064283

$ws.Range("D10").Value = 0.1207336743123858
$ws.Range("E10").Value = -0.01629201897298416

$ws.Range("D11").Value = 0.2497125447989407
$ws.Range("E11").Value = -0.007577910391209586

$ws.Range("D12").Value = 0.1063731146100942
$ws.Range("E12").Value = -0.01110892646480532

$ws.Range("E13").Value = -0.01165382476013865

$ws.Protect()
